$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: externalId values (ceo / emp1 / emp2) ---
$ws.Cells.Item(2,1).Value = "ceo"
$ws.Cells.Item(3,1).Value = "emp1"
$ws.Cells.Item(4,1).Value = "emp2"

# --- Column C: ssn values ---
$ws.Cells.Item(2,3).Value = "090977-954P"
$ws.Cells.Item(3,3).Value = "161165-951M"
$ws.Cells.Item(4,3).Value = "110674-9046"

# --- Columns D & E: callName / lastName, row by row ---
$ws.Cells.Item(2,4).Value = "Cecily"
$ws.Cells.Item(2,5).Value = "Ceo"
$ws.Cells.Item(3,4).Value = "Adam"
$ws.Cells.Item(3,5).Value = "Ant"
$ws.Cells.Item(4,4).Value = "Betty"
$ws.Cells.Item(4,5).Value = "Boo"

# --- Column F: emailAddress (with a fill style touch, like the source workbook) ---
$ws.Cells.Item(2,6).Value = "ceo@company.com"
$ws.Cells.Item(2,6).Interior.ColorIndex = -4142
$ws.Cells.Item(3,6).Value = "adam.ant@company.com"
$ws.Cells.Item(3,6).Interior.ColorIndex = -4142
$ws.Cells.Item(4,6).Value = "betty.boo@company.com"
$ws.Cells.Item(4,6).Interior.ColorIndex = -4142

# --- Column I: localPhoneNumber, stored as text (leading apostrophe => quote-prefix style) ---
$ws.Cells.Item(2,9).Value = "'0101234567"
$ws.Cells.Item(3,9).Value = "'0101122334"
$ws.Cells.Item(4,9).Value = "'0107654321"

# --- Row 2: startDate / endDate / department / departmentStart ---
$ws.Cells.Item(2,11).Value = 42401
$ws.Cells.Item(2,11).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2,12).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2,13).Value = "dep1"
$ws.Cells.Item(2,14).Value = 42736
$ws.Cells.Item(2,14).NumberFormat = "mm-dd-yy"

# --- Row 3: startDate / endDate / department / departmentStart / supervisor / supervisorStart ---
$ws.Cells.Item(3,11).Value = 43102
$ws.Cells.Item(3,11).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(3,12).Value = 44196
$ws.Cells.Item(3,12).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(3,13).Value = "dep1"
$ws.Cells.Item(3,14).Value = 43132
$ws.Cells.Item(3,14).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(3,15).Value = "ceo"
$ws.Cells.Item(3,16).Value = 43132
$ws.Cells.Item(3,16).NumberFormat = "mm-dd-yy"

# --- Row 4: startDate / endDate / department / departmentStart / supervisor / supervisorStart ---
$ws.Cells.Item(4,11).Value = 42767
$ws.Cells.Item(4,11).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(4,12).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(4,13).Value = "dep2"
$ws.Cells.Item(4,14).Value = 42767
$ws.Cells.Item(4,14).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(4,15).Value = "ceo"
$ws.Cells.Item(4,16).Value = 42767
$ws.Cells.Item(4,16).NumberFormat = "mm-dd-yy"

# --- Column widths (bestFit-like custom widths for the email + department columns) ---
$ws.Columns.Item(6).ColumnWidth = 21.666666666666668
$ws.Columns.Item(13).ColumnWidth = 10

# --- Sheet view / selection: land on A4, no frozen top-left override ---
$ws.Range("A4").Select()
